# Build site at 2022-09-26 16:07:08 UTC
# LOQ4004.xlsx content update: the Objetivos/Programa resumido/Programa/
# Avaliacao block was rewritten and one trailing row was removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 10 ("Objetivos:") keeps its place but its B/C text changes to
# the docente responsavel string.
# ------------------------------------------------------------------
$ws.Range("B10").Value = "8643537 - Fabio Rodolfo Miguel Batista"
$ws.Range("C10").Value = "8643537 - Fabio Rodolfo Miguel Batista"

# ------------------------------------------------------------------
# Everything from row 13 down gets reshuffled/rewritten and the sheet
# ends up one row shorter (old row 26 disappears). The cleanest way to
# reach the exact target cell layout (and avoid leaving stray empty
# <c> elements behind) is to delete rows 13-26 completely and then
# populate fresh rows 13-25 from scratch.
# ------------------------------------------------------------------
$ws.Rows("13:26").Delete()

function Set-TextValue {
    # Writes $Text into $Range as a plain text value even when it looks
    # like a date (e.g. "01/01/2013"), by staging it through a scratch
    # cell that is forced to text format first, then copying just the
    # resulting value in - this avoids Excel's automatic date parsing.
    param($Range, [string]$Text)

    $scratch = $ws.Range("ZZ1000")
    $scratch.NumberFormat = "@"
    $scratch.Value = $Text
    $scratch.Copy()
    $Range.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $scratch.Delete()
}

function Set-LabelRow {
    param([int]$Row, $A, $B, $C, $Height)

    if ($A -ne $null) {
        Set-TextValue $ws.Range("A$Row") $A
        $ws.Range("A10").Copy()
        $ws.Range("A$Row").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    }
    if ($B -ne $null) {
        Set-TextValue $ws.Range("B$Row") $B
        $ws.Range("B10").Copy()
        $ws.Range("B$Row").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    }
    if ($C -ne $null) {
        Set-TextValue $ws.Range("C$Row") $C
        $ws.Range("C10").Copy()
        $ws.Range("C$Row").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    }
    if ($Height -ne $null) {
        $ws.Rows($Row).RowHeight = $Height
    }
}

Set-LabelRow 13 "Programa resumido:" "Semestral" "Semestral" 60
Set-LabelRow 14 "Short syllabus:" $null $null 60
Set-LabelRow 15 "Programa:" "01/01/2013" "01/01/2013" 120
Set-LabelRow 16 "Syllabus:" $null $null 120
Set-LabelRow 17 "Avaliação:" $null $null $null
Set-LabelRow 18 "Método:" "8643537 - Fabio Rodolfo Miguel Batista" "8643537 - Fabio Rodolfo Miguel Batista" 60
Set-LabelRow 19 "Critério:" "Duas provas escritas: P1 e P2" "Duas provas escritas: P1 e P2" 60
Set-LabelRow 20 "Norma de recuperação:" "Média das notas obtidas nas duas provas: N1=(P1 + P2)/2" "Média das notas obtidas nas duas provas: N1=(P1 + P2)/2" 60
Set-LabelRow 21 "Bibliografia:" "Uma prova escrita: REC`nMédia das notas N1 e REC:N2=(N1+REC)/2" "Uma prova escrita: REC`nMédia das notas N1 e REC:N2=(N1+REC)/2" 120
Set-LabelRow 22 "Requisitos:" $null $null $null
Set-LabelRow 23 $null "LOB1006 -  Cálculo IV  (Requisito fraco)`n" "LOB1006 -  Cálculo IV  (Requisito fraco)`n" 30
Set-LabelRow 24 $null "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n" "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n" 30
Set-LabelRow 25 $null "LOQ4009 -  Instrumentação na Industria Química  (Requisito fraco)`n" "LOQ4009 -  Instrumentação na Industria Química  (Requisito fraco)`n" 30
